$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new date column (D) for 2017-11-16, matching the existing date in C2 (2017-11-15)
$ws.Range("D2").Value = 43055
$ws.Range("D2").NumberFormat = "d-mmm"

# Fill in hours logged for that new date for Rick (row 3), Stijn (row 5), Stan (row 7)
$ws.Range("D3").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("D7").Value = 3

# Update the active selection/cell
$ws.Range("E9").Select()
